# Insert a new row at position 191, pushing existing rows 191-218 down to 192-219.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new record.
$row = 191

$ws.Cells.Item($row, 1).Value  = 3
$ws.Cells.Item($row, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value  = "Coquimbo"
$ws.Cells.Item($row, 4).Value  = 44491
$ws.Cells.Item($row, 5).Value  = 5
$ws.Cells.Item($row, 6).Value  = 100112043
$ws.Cells.Item($row, 7).Value  = "Pepino ensalada"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 130
$ws.Cells.Item($row, 11).Value = 7000
$ws.Cells.Item($row, 12).Value = 8000
$ws.Cells.Item($row, 13).Value = 7538
$ws.Cells.Item($row, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 108
$ws.Cells.Item($row, 17).Value = 70
$ws.Cells.Item($row, 18).Value = "Hortaliza"
